# Apply the cryptos.xlsx data refresh described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the Price (D) cells that parse as plain decimals as Text first,
# so Excel does not silently convert them to numbers (the source data
# keeps these as literal strings, matching the original workbook).
$priceCellsNeedingTextFormat = @(5,6,9,10,14,18,19,20,21,22,27,30,31,32,34,35,36,37,39,41,43,44,49,50)
foreach ($r in $priceCellsNeedingTextFormat) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range('D2').Value = '62.928.69'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '3.454.70'
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '580.74'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').Value = '150.68'
$ws.Range('E6').Value = '  +2.40%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  +1.49%  '
$ws.Range('D9').Value = '8.08'
$ws.Range('E9').Value = '  +6.39%  '
$ws.Range('D10').Value = '0.125'
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('E11').Value = '  +3.67%  '
$ws.Range('D12').Value = '4.044.54'
$ws.Range('E12').Value = '  -0.39%  '
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('D14').Value = '28.37'
$ws.Range('E14').Value = '  -4.88%  '
$ws.Range('D15').Value = '3.459.17'
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('E16').Value = '  +1.68%  '
$ws.Range('D17').Value = '62.936.47'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('D18').Value = '6.44'
$ws.Range('E18').Value = '  +1.64%  '
$ws.Range('D19').Value = '14.61'
$ws.Range('E19').Value = '  +1.83%  '
$ws.Range('D20').Value = '9.03'
$ws.Range('E20').Value = '  -2.17%  '
$ws.Range('D21').Value = '388.57'
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').Value = '0.570'
$ws.Range('E22').Value = '  +1.71%  '
$ws.Range('E23').Value = '  +0.87%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = '3.590.29'
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('D27').Value = '0.186'
$ws.Range('E27').Value = '  +4.30%  '
$ws.Range('E28').Value = '  +2.16%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').Value = '8.03'
$ws.Range('E30').Value = '  -1.63%  '
$ws.Range('D31').Value = '2.14'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('E33').Value = '  -2.32%  '
$ws.Range('D34').Value = '23.35'
$ws.Range('E34').Value = '  -1.58%  '
$ws.Range('D35').Value = '5.44'
$ws.Range('E35').Value = '  +2.89%  '
$ws.Range('D36').Value = '1.65'
$ws.Range('E36').Value = '  +4.15%  '
$ws.Range('D37').Value = '32.16'
$ws.Range('E37').Value = '  +2.29%  '
$ws.Range('E38').Value = '  -1.75%  '
$ws.Range('D39').Value = '169.12'
$ws.Range('E39').Value = '  -0.19%  '
$ws.Range('D40').Value = '3.488.17'
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('D41').Value = '0.0782'
$ws.Range('E41').Value = '  +2.28%  '
$ws.Range('E42').Value = '  -1.25%  '
$ws.Range('D43').Value = '42.78'
$ws.Range('E43').Value = '  +1.04%  '
$ws.Range('D44').Value = '1.71'
$ws.Range('E44').Value = '  -0.50%  '
$ws.Range('E45').Value = '  -1.58%  '
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('D47').Value = '2.563.65'
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('E48').Value = '  +2.97%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '22.92'
$ws.Range('E49').Value = '  -0.87%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').Value = '2.23'
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('E51').Value = '  -0.06%  '
